$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet: append new Name/Value rows (13-19) describing the new
# business-exception / set-exception config codes plus a couple of new
# path/email settings used by the new sequences
# (move_VendorListExcel, Obtain_data, Reset_test_data, Write_Data).
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")

$wsSettings.Range("A13").Value = "BECode_Credential"
$wsSettings.Range("B13").Value = "BE#001"

$wsSettings.Range("A18").Value = "Client_Email"
$wsSettings.Range("B18").Value = "sarashattra@gmail.com"

$wsSettings.Range("B14").Value = "BE#003"
$wsSettings.Range("A14").Value = "BECode_Amount"

$wsSettings.Range("A16").Value = "BECode_NoInvoices"
$wsSettings.Range("B16").Value = "BE#005"

$wsSettings.Range("A15").Value = "BECode_NoEURInvoices"
$wsSettings.Range("B15").Value = "BE#004"

$wsSettings.Range("A17").Value = "SECode_ResetTestData"
$wsSettings.Range("B17").Value = "SE#006"

$wsSettings.Range("A19").Value = "Path_Output_Vendors"
$wsSettings.Range("B19").Value = "C:\Output_Vendors"

# ---------------------------------------------------------------------
# Constants sheet: bump MaxConsecutiveSystemExceptions from 0 to 4.
# ---------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Value = 4

# Leave each sheet's selection the way the author left it, then return
# to (and leave active) the Settings tab, matching the saved workbook.
$wsConstants.Activate() | Out-Null
$wsConstants.Range("B13").Select() | Out-Null

$wsSettings.Activate() | Out-Null
$wsSettings.Range("A9").Select() | Out-Null
